$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Vietnamese headers (row 1) -- order chosen to match shared-string append order
$ws.Range("A1").Value = "Mã đề tài(*)"
$ws.Range("B1").Value = "Tiêu đề(*)"
$ws.Range("D1").Value = "Mã đợt"
$ws.Range("E1").Value = "Mã GVHD"
$ws.Range("C1").Value = "Số lượng SV(*)"
$ws.Range("F1").Value = "Mã GVPB"
$ws.Range("G1").Value = "Mô tả"

# Remove old Description column header (H) entirely
$ws.Range("H1").Clear()

# Re-order row 2 sample data to match new columns
$ws.Range("A2").Value = "2023-KL-001"
$ws.Range("B2").Value = "Xây dựng hệ thống nhận diện chó "
$ws.Range("C2").Value = "2"

# D column used to carry a date format (Thesis Defense Date) - drop that,
# it's now a plain text column like the others
$ws.Range("D1").NumberFormat = "@"

# Set font size to 12 (matches new default font, size 12) - only the
# populated cells so we don't materialize empty styled cells
$ws.Range("A1:G1").Font.Size = 12
$ws.Range("A2:C2").Font.Size = 12

# "So luong SV" (count) column is centered
$ws.Range("C2").HorizontalAlignment = -4108

# Column widths
$ws.Range("A:A").ColumnWidth = 16.109375
$ws.Range("B:B").ColumnWidth = 47.33203125
$ws.Range("C:C").ColumnWidth = 18.33203125
$ws.Range("D:F").ColumnWidth = 18.77734375
$ws.Range("G:G").ColumnWidth = 61.33203125

$wb.Save()
